# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values computed/regenerated for rows 2-17
$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 3
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
